$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '68.607.20'
$ws.Range("E2").Value = '  +0.97%  '
$ws.Range("D3").Value = '3.719.20'
$ws.Range("E3").Value = '  -2.61%  '
$ws.Range("E4").Value = '  +0.11%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '601.37'
$ws.Range("E5").Value = '  -0.20%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '166.79'
$ws.Range("E6").Value = '  -4.82%  '
$ws.Range("D7").Value = '3.717.45'
$ws.Range("E7").Value = '  -2.79%  '
$ws.Range("E8").Value = '  +0.04%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.535'
$ws.Range("E9").Value = '  +0.98%  '
$ws.Range("E10").Value = '  +2.79%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.37'
$ws.Range("E11").Value = '  +1.69%  '
$ws.Range("E12").Value = '  -1.87%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '37.98'
$ws.Range("E13").Value = '  -1.94%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000245'
$ws.Range("E14").Value = '  -0.52%  '
$ws.Range("D15").Value = '4.345.29'
$ws.Range("E15").Value = '  -2.25%  '
$ws.Range("D16").Value = '3.732.02'
$ws.Range("E16").Value = '  -2.05%  '
$ws.Range("D17").Value = '68.564.01'
$ws.Range("E17").Value = '  +0.86%  '
$ws.Range("E18").Value = '  -0.69%  '
$ws.Range("E19").Value = '  +0.13%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.10'
$ws.Range("E20").Value = '  +1.21%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '496.19'
$ws.Range("E21").Value = '  +0.47%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '10.34'
$ws.Range("E22").Value = '  +12.15%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.725'
$ws.Range("E23").Value = '  -3.89%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '84.97'
$ws.Range("E24").Value = '  -1.50%  '
$ws.Range("E25").Value = '  -4.45%  '
$ws.Range("E26").Value = '  -2.78%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '12.45'
$ws.Range("E27").Value = '  +0.38%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.09'
$ws.Range("E28").Value = '  -1.80%  '
$ws.Range("E29").Value = '  -0.17%  '
$ws.Range("B30").Value = 'ImmutableX'
$ws.Range("C30").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.56'
$ws.Range("E30").Value = '  +4.12%  '
$ws.Range("B31").Value = 'PancakeSwap'
$ws.Range("C31").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.95'
$ws.Range("E31").Value = '  -1.21%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '7.90'
$ws.Range("E32").Value = '  +1.40%  '
$ws.Range("E33").Value = '  -5.91%  '
$ws.Range("D34").Value = '3.866.72'
$ws.Range("E34").Value = '  -2.03%  '
$ws.Range("E35").Value = '  -1.78%  '
$ws.Range("D36").Value = '3.655.13'
$ws.Range("E36").Value = '  -2.39%  '
$ws.Range("E37").Value = '  +0.13%  '
$ws.Range("E38").Value = '  -0.61%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.83'
$ws.Range("E39").Value = '  -0.73%  '
$ws.Range("E40").Value = '  -1.69%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.324'
$ws.Range("E41").Value = '  -1.81%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '435.50'
$ws.Range("E42").Value = '  -4.66%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '49.05'
$ws.Range("E43").Value = '  -0.34%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.99'
$ws.Range("E44").Value = '  -0.88%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.88'
$ws.Range("E45").Value = '  -0.52%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '8.48'
$ws.Range("E46").Value = '  +0.51%  '
$ws.Range("B47").Value = 'USDe'
$ws.Range("C47").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.00'
$ws.Range("E47").Value = '  +0.03%  '
$ws.Range("B48").Value = 'Arweave'
$ws.Range("C48").Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '40.69'
$ws.Range("E48").Value = '  -0.60%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '141.56'
$ws.Range("E49").Value = '  +2.09%  '
$ws.Range("E50").Value = '  +0.01%  '
$ws.Range("D51").Value = '2.745.91'
$ws.Range("E51").Value = '  -3.91%  '
